$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.884.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.639.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5049'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2571'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06400'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07776'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.285'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.648.04'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5438'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0₅7874'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.952.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '197.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.396'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.959'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.030'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.006'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.867'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '140.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1143'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.862'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.238'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05012'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.257'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.193'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.363'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.8939'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.134.01'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5529'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01555'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.005'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.681'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8144'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("E44").Value = '  +13.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.777.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4534'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.17'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05084'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.007'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09533'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.05%  '
